# Commit enne kodutöö 3 alustamist
#
# Inserts a new "Nädal 3" worksheet at the front of the workbook (a fresh
# copy of the "Nädal 2" template with its logged time entries cleared out),
# and fills in the last logged entry (#11) that was previously missing on
# "Nädal 2".

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "Nädal 3" sheet ------------------------------------
# Copy "Nädal 2" (keeps all formatting/styles/merged cells/formulas) and
# place the copy immediately before it, then rename it.
$wsWeek2src = $wb.Worksheets.Item("Nädal 2")
$wsWeek2src.Copy($wsWeek2src)
$wsWeek3 = $wb.Worksheets.Item(1)
$wsWeek3.Name = "Nädal 3"

# The COM object for the original sheet doesn't track the rename of the
# *other* (copied) sheet reliably, so re-resolve "Nädal 2" by name now that
# the tab layout has settled.
$wsWeek2 = $wb.Worksheets.Item("Nädal 2")
$wsWeek1 = $wb.Worksheets.Item("Nädal 1")

# Clear out the copied time-log rows (7-16) so the new week starts blank,
# but keep the row index column (A) and the pre-formatted blank columns.
$wsWeek3.Range("B7:D16").ClearContents()
$wsWeek3.Range("F7:G16").ClearContents()

# New week's start date.
$wsWeek3.Range("G4").Value = 43877

# Tidy up the view: no frozen/scrolled top-left cell, selection on the
# (still blank) log rows.
$wsWeek3.Range("A5:J5").Select()

# --- 2. Fill in the missing last entry on "Nädal 2" ------------------------
$wsWeek2.Range("G4").Value = 43870

$wsWeek2.Range("B17").Value = 43870
$wsWeek2.Range("C17").Value = 0.78819444444444453
$wsWeek2.Range("D17").Value = 0.80555555555555547
$wsWeek2.Range("F17").Value = 25
$wsWeek2.Range("G17").Value = "Razor page, MVC mõisted"

$wsWeek2.Range("G4:J4").Select()

# --- 3. Cosmetic selection update on "Nädal 1" ------------------------------
$wsWeek1.Range("G4:J4").Select()

# Leave focus on the newly added sheet, matching the commit's "about to
# start homework 3" checkpoint.
$wsWeek3.Select()
